# Add daily power records
# Fills in "Start Time" (B) / "End Time" (C) entries for the days that were
# previously blank (rows 106-111 of the comforter-cda sheet): five full
# days with no recorded activity (B=0, C=0) followed by a day where work
# has only just started (B set, C still empty -> shared formulas in
# D/E/F recompute automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 106-110: both a start and an end time of 0 (no activity recorded yet)
106..110 | ForEach-Object {
    $row = $_
    $ws.Cells.Item($row, 2).Value = 0   # column B - Start Time
    $ws.Cells.Item($row, 3).Value = 0   # column C - End Time
}

# Row 111: only a start time has been recorded so far
$ws.Range("B111").Value = 0.77430555555555547

# Match the workbook's new selection/scroll position
$ws.Range("D111").Select()
